$d = $word.ActiveDocument

# Merge the split "<id>...</id>" runs back into a single run for each
# downloaded tc/tcn/tl entry (p009v_1 and p009v_2). Word's Find/Replace
# collapses the matched range into one run, adopting the formatting of
# the first character of the match (Courier New / 7f6000 / sz 18, the
# <id> tag's own formatting), which is exactly the target state.
$d.Content.Find.Execute("<id>p009v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p009v_1</id>", 2)
$d.Content.Find.Execute("<id>p009v_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p009v_2</id>", 2)
